$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 15614.315
$ws.Range("I19").Value = 2210.875
$ws.Range("J19").Value = 25362.273
$ws.Range("K19").Value = 2210.875
$ws.Range("L19").Value = 25362.273
$ws.Range("M19").Value = -2035.875
$ws.Range("N19").Value = -25712.273

$ws.Range("H100").Value = 6529.8
$ws.Range("I100").Value = 6699.778
$ws.Range("K100").Value = 6699.778
$ws.Range("M100").Value = -6158.778

$ws.Range("H101").Value = 296.0625
$ws.Range("J101").Value = 145
$ws.Range("L101").Value = 435
$ws.Range("N101").Value = -3679

$ws.Range("H103").Value = 473.58334
$ws.Range("I103").Value = 303.66666
$ws.Range("J103").Value = 983.3333
$ws.Range("K103").Value = 910.9999799999999
$ws.Range("L103").Value = 2949.9999
$ws.Range("M103").Value = -324.9999799999999
$ws.Range("N103").Value = -4121.9999

$ws.Range("H111").Value = 42310.727
$ws.Range("I111").Value = 3549.3333
$ws.Range("K111").Value = 10647.9999
$ws.Range("M111").Value = -7580.999899999999

$ws.Range("H129").Value = 1312
$ws.Range("I129").Value = 1175.7778
$ws.Range("K129").Value = 3527.3334
$ws.Range("M129").Value = 1472.6666

$ws.Range("H132").Value = 4644.9033
$ws.Range("I132").Value = 3573.6296
$ws.Range("K132").Value = 10720.8888
$ws.Range("M132").Value = -8190.888800000001

$ws.Range("H137").Value = 2709.4707
$ws.Range("I137").Value = 2279.24
$ws.Range("K137").Value = 6837.719999999999
$ws.Range("M137").Value = -4287.719999999999

$ws.Range("H141").Value = 7865.5713
$ws.Range("I141").Value = 7865.5713
$ws.Range("K141").Value = 23596.7139
$ws.Range("M141").Value = -18416.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3329.8
$ws.Range("I61").Value = 1971.2858
$ws.Range("K61").Value = 1971.2858
$ws.Range("M61").Value = -1759.2858

$ws.Range("H74").Value = 66670940
$ws.Range("I74").Value = 83336550
$ws.Range("K74").Value = 83336550
$ws.Range("M74").Value = -83335676

$ws.Range("H77").Value = 66670940
$ws.Range("I77").Value = 83336550
$ws.Range("K77").Value = 416682750
$ws.Range("M77").Value = -416678382

$ws.Range("H122").Value = 2006.674
$ws.Range("I122").Value = 1200.2646
$ws.Range("K122").Value = 3600.7938
$ws.Range("M122").Value = -1150.7938

$ws.Range("H127").Value = 93593.42999999999
$ws.Range("I127").Value = 37999.4
$ws.Range("K127").Value = 37999.4
$ws.Range("M127").Value = -33039.4

$ws.Range("H132").Value = 2674.48
$ws.Range("I132").Value = 1625
$ws.Range("J132").Value = 4904.625
$ws.Range("K132").Value = 4875
$ws.Range("L132").Value = 14713.875
$ws.Range("M132").Value = -2345
$ws.Range("N132").Value = -19773.875

$ws.Range("H136").Value = 3329.8
$ws.Range("I136").Value = 1971.2858
$ws.Range("K136").Value = 5913.857400000001
$ws.Range("M136").Value = -3363.857400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 8000
$ws.Range("I26").Value = 8000
$ws.Range("K26").Value = 8000
$ws.Range("M26").Value = -7708

$ws.Range("H42").Value = 299684
$ws.Range("J42").Value = 299684
$ws.Range("L42").Value = 299684
$ws.Range("N42").Value = -300340

$ws.Range("H86").Value = 1746.3889
$ws.Range("I86").Value = 1303.625
$ws.Range("K86").Value = 1303.625
$ws.Range("M86").Value = -180.625

$ws.Range("H89").Value = 1746.3889
$ws.Range("I89").Value = 1303.625
$ws.Range("K89").Value = 6518.125
$ws.Range("M89").Value = -902.125

$ws.Range("H105").Value = 2119.8
$ws.Range("I105").Value = 1885.05
$ws.Range("J105").Value = 2589.3
$ws.Range("K105").Value = 1885.05
$ws.Range("L105").Value = 2589.3
$ws.Range("M105").Value = -138.05
$ws.Range("N105").Value = -6083.3

$ws.Range("H134").Value = 4250.8125
$ws.Range("I134").Value = 3944.5
$ws.Range("K134").Value = 11833.5
$ws.Range("M134").Value = -9298.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 6750
$ws.Range("J2").Value = 6750
$ws.Range("L2").Value = 6750
$ws.Range("N2").Value = -6976

$ws.Range("H31").Value = 3269.48
$ws.Range("I31").Value = 2534.3125
$ws.Range("J31").Value = 3468.8474
$ws.Range("K31").Value = 2534.3125
$ws.Range("L31").Value = 3468.8474
$ws.Range("M31").Value = -2239.3125
$ws.Range("N31").Value = -4058.8474

$ws.Range("H34").Value = 3269.48
$ws.Range("I34").Value = 2534.3125
$ws.Range("J34").Value = 3468.8474
$ws.Range("K34").Value = 2534.3125
$ws.Range("L34").Value = 3468.8474
$ws.Range("M34").Value = -2332.3125
$ws.Range("N34").Value = -3872.8474

$ws.Range("H58").Value = 4037.4707
$ws.Range("I58").Value = 1945.9
$ws.Range("K58").Value = 1945.9
$ws.Range("M58").Value = -1742.9

$ws.Range("H99").Value = 8742.23
$ws.Range("I99").Value = 6864.9
$ws.Range("K99").Value = 6864.9
$ws.Range("M99").Value = -5366.9

$ws.Range("H122").Value = 2622.7144
$ws.Range("J122").Value = 1282
$ws.Range("L122").Value = 3846
$ws.Range("N122").Value = -8746

$ws.Range("H126").Value = 8742.23
$ws.Range("I126").Value = 6864.9
$ws.Range("K126").Value = 20594.7
$ws.Range("M126").Value = -18124.7

$ws.Range("H132").Value = 3570.2
$ws.Range("I132").Value = 2814.7778
$ws.Range("K132").Value = 8444.3334
$ws.Range("M132").Value = -5914.3334

$ws.Range("H134").Value = 3658.4
$ws.Range("I134").Value = 2598.4707
$ws.Range("K134").Value = 7795.4121
$ws.Range("M134").Value = -5260.4121

$ws.Range("H135").Value = 66311.14
$ws.Range("J135").Value = 66311.14
$ws.Range("L135").Value = 66311.14
$ws.Range("N135").Value = -76451.14

$ws.Range("H136").Value = 4037.4707
$ws.Range("I136").Value = 1945.9
$ws.Range("K136").Value = 5837.700000000001
$ws.Range("M136").Value = -3287.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1250.1765
$ws.Range("J2").Value = 1002
$ws.Range("L2").Value = 6012
$ws.Range("N2").Value = -6238

$ws.Range("H51").Value = 1244.6666
$ws.Range("I51").Value = 1244.6666
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 3733.9998
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -3273.9998
$ws.Range("N51").ClearContents()

$ws.Range("H134").Value = 4290.6
$ws.Range("J134").Value = 9796.666999999999
$ws.Range("L134").Value = 29390.001
$ws.Range("N134").Value = -39530.001

$ws.Range("H138").Value = 2503316.5
$ws.Range("I138").Value = 7504250
$ws.Range("K138").Value = 22512750
$ws.Range("M138").Value = -22507610

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 649.4286
$ws.Range("J2").Value = 676.6667
$ws.Range("L2").Value = 676.6667
$ws.Range("N2").Value = -902.6667

$ws.Range("H122").Value = 1738.2903
$ws.Range("I122").Value = 1460.4783
$ws.Range("J122").Value = 2537
$ws.Range("K122").Value = 4381.4349
$ws.Range("L122").Value = 7611
$ws.Range("M122").Value = -1931.4349
$ws.Range("N122").Value = -12511

$ws.Range("H132").Value = 3546.1538
$ws.Range("I132").Value = 3429.7
$ws.Range("K132").Value = 10289.1
$ws.Range("M132").Value = -7759.099999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 17860204
$ws.Range("I7").Value = 27780292
$ws.Range("J7").Value = 4045.7
$ws.Range("K7").Value = 27780292
$ws.Range("L7").Value = 4045.7
$ws.Range("M7").Value = -27780180
$ws.Range("N7").Value = -4269.7

$ws.Range("H26").Value = 15633
$ws.Range("I26").Value = 15633
$ws.Range("K26").Value = 15633
$ws.Range("M26").Value = -15338

$ws.Range("H126").Value = 17860204
$ws.Range("I126").Value = 27780292
$ws.Range("J126").Value = 4045.7
$ws.Range("K126").Value = 83340876
$ws.Range("L126").Value = 12137.1
$ws.Range("M126").Value = -83338406
$ws.Range("N126").Value = -17077.1

$ws.Range("H132").Value = 22732054
$ws.Range("I132").Value = 30305648
$ws.Range("K132").Value = 90916944
$ws.Range("M132").Value = -90914414

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 6583.2856
$ws.Range("I113").Value = 7670.7144
$ws.Range("K113").Value = 23012.1432
$ws.Range("M113").Value = -20842.1432

$ws.Range("H122").Value = 1777.0857
$ws.Range("I122").Value = 1533.2963
$ws.Range("J122").Value = 2599.875
$ws.Range("K122").Value = 4599.8889
$ws.Range("L122").Value = 7799.625
$ws.Range("M122").Value = -2149.8889
$ws.Range("N122").Value = -12699.625

$ws.Range("H132").Value = 5011.436
$ws.Range("I132").Value = 4749.909
$ws.Range("K132").Value = 14249.727
$ws.Range("M132").Value = -11719.727

$ws.Range("H136").Value = 4750.057
$ws.Range("J136").Value = 4655.7144
$ws.Range("L136").Value = 13967.1432
$ws.Range("N136").Value = -19067.1432
